$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Text fixes (typos / grammar / spelling Americanization) that the
#    commit message describes as "Typo fixed, Grammerly used". Each
#    paragraph body is replaced in one shot so the corrected text ends
#    up back in a single run.
# ------------------------------------------------------------------

# Paragraph: "present time..."
$d.Content.Find.Execute("At the present time, robotics grows day by day and robots have a new feature at every turn. For example, robots walk like human instead of moving by using wheel as the past. Also, one of the important features is ‘sensing environment’ for robots. It is wanted that robots’ sense like human and comprehend their domain or surroundings. ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "At the present time, robotics grows day by day and robots have a new feature at every turn. For example, robots walk like human instead of moving by using the wheel as the past. Also, one of the important features is ‘sensing environment’ for robots. It is wanted that robots’ sense like human and comprehend their domain or surroundings.", 2) | Out-Null

# Paragraph: "this project..."
$d.Content.Find.Execute("At this project, sensing environment is emphasized.  A device or robot is wanted to sense its environment and extracts the plan of its surrounding from this sense. The device rambles on the region that is determined by standards and creates a plan of the region. The plan is sent the computer screen. In addition, İt is actually a race. Two devices extract plan of the same region and the aims are better plan and minimum time.  If it needs to talk about the environment, A region is closed with twelve walls. The walls are fifty centimetre long. The region includes 8 objects with the shape of cylinders with ten- and five-centimetre diameter, square prism with seven-centimetre edge length and prism with triangular base of eight-centimetre edge length. In other respect, the robot is communicated with computer that shows the plan at one direction. There is no information flow from computer to robot. İt means that the robot moves with autonomous to extract plan. ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "At this project, sensing environment is emphasized.  A device or robot is wanted to sense its environment and extracts the plan of its surrounding from this sense. The device rambles on the region that is determined by standards and creates a plan of the region. The plan is sent the computer screen. In addition, it is actually a race. Two devices extract plan of the same region and the aims are the better plan and minimum time.  If it needs to talk about the environment, A region is closed with twelve walls. The walls are fifty centimeters long. The region includes 8 objects with the shape of cylinders with ten- and five-centimeter diameter, square prisms with seven-centimeter edge length and prism with a triangular base of eight-centimeter edge length. In other respect, the robot is communicated with the computer that shows the plan in one direction. There is no information flow from the computer to the robot. It means that the robot moves with autonomous to extract plan.", 2) | Out-Null

# Paragraph: "Moreover..."
$d.Content.Find.Execute("Moreover, the mapping is required to some analytical and mathematical algorithm. The robot’s rambles on region while sensor or camera is taking a information about the environment. So, there are two data flow from robots. One of them is position information about robot and other is surrounding distance information. These data are mixed with algorithm to extract plan.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Moreover, the mapping is required to some analytical and mathematical algorithm. The robot’s rambles on the region while the sensor or camera is taking an information about the environment. So, there are two data flow from robots. One of them is position information about robot and other is surrounding distance information. These data are mixed with the algorithm to extract plan.", 2) | Out-Null

# Paragraph: "Finally..."
$d.Content.Find.Execute("Finally, the project gets up to date with development of robotics and AI. Also, it needs to base of some algebraic transformation alongside optimization on software. ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Finally, the project gets up to date with the development of robotics and AI. Also, it needs to base of some algebraic transformation alongside optimization on software.", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of the last paragraph
#    (the previous last-edit location) to mid-way through the first
#    body paragraph (the new last-edit location), matching the XML
#    diff.
# ------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$ip = $d.Content
$ip.Find.Execute("robotics grows da", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$ip.Collapse(0)
$d.Bookmarks.Add("_GoBack", $ip)
